# Daily attendance processing - reorder the "Recorded By" (column G) names.
# For every data row, the comma-separated list of recorders is rotated right
# by one position: the last name in the list is moved to the front, and the
# rest keep their existing relative order. Single-value cells are unaffected
# (rotating a 1-item list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $old = $cell.Text

    if ([string]::IsNullOrEmpty($old)) { continue }

    $parts = $old.Split(",")
    if ($parts.Length -lt 2) { continue }

    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $lastItem = $trimmed[$trimmed.Length - 1]
    $rest = $trimmed[0..($trimmed.Length - 2)]
    $newParts = @($lastItem) + $rest
    $newVal = [string]::Join(", ", $newParts)

    $cell.Value = $newVal
}
